# UGN-356 - add "date" column (dateFormat prop test fixture)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell: new "date" column, styled with the (new) default Calibri/black font
$ws.Range("C1").Value = "date"
$ws.Range("C1").Font.Color = 0

# Data cells: serials for 2020-03-03 / 2010-04-04 / 1994-02-27, date-formatted
$ws.Range("C2").Value = 43893
$ws.Range("C3").Value = 40272
$ws.Range("C4").Value = 34392

$ws.Range("C2").Font.Color = 0
$ws.Range("C2").NumberFormat = "mm-dd-yy"

# Re-use the exact same style record for C3:C4 (copy/paste-format instead of
# re-assigning Font/NumberFormat, which would otherwise mint a fresh cellXfs
# entry per cell)
$ws.Range("C2").Copy()
$ws.Range("C3:C4").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Move the active selection to the last populated cell
[void]$ws.Range("C4").Select()
